$d = $word.ActiveDocument

# Locate the paragraph that ends with the "Hierarchies: metaclass ..." sentence;
# the new content is inserted right after it (and before the existing blank
# paragraph that precedes "Model Application domains ...").
$target = "Hierarchies: metaclass / role, class, instance, occurrence (parent, children, previous, next, attribute, value). Encoding. Comparisons. Functional traversal (streams)."

$rng = $d.Content
$found = $rng.Find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $rng.Collapse(0)
    $pos = $rng.Start

    # Six new paragraphs: blank line, sentence, blank line, statement, blank
    # line, statement - matching the authored edit exactly.
    $lines = @(
        "",
        "Dimensional statements. Occurrences contexts: Events. Dimensional contexts: Context occurrences hierarchy for order relations assertions.",
        "",
        "(Mapping, Kind, Role, Statement);",
        "",
        "(Mapping, Mapping super / parent / dimension, Kind unit, Role measure);"
    )

    foreach ($line in $lines) {
        $curRng = $d.Range($pos, $pos)
        $curRng.InsertParagraphAfter()
        $pos = $pos + 1
        if ($line -ne "") {
            $textRng = $d.Range($pos, $pos)
            $textRng.InsertAfter($line)
            $pos = $pos + $line.Length
        }
    }
}

Write-Output "paragraphs: $($d.Paragraphs.Count)"
